$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.997.67'
$ws.Range("E2").Value = '  +2.92%  '
$ws.Range("D3").Value = '3.050.24'
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '519.64'
$ws.Range("E5").Value = '  +3.46%  '
$ws.Range("D6").Value = '141.83'
$ws.Range("E6").Value = '  +5.55%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +4.11%  '
$ws.Range("D9").Value = '7.51'
$ws.Range("E9").Value = '  +2.74%  '
$ws.Range("E10").Value = '  +5.84%  '
$ws.Range("E11").Value = '  +5.56%  '
$ws.Range("D12").Value = '3.576.03'
$ws.Range("E12").Value = '  +2.62%  '
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("D14").Value = '26.73'
$ws.Range("E14").Value = '  +6.77%  '
$ws.Range("E15").Value = '  +13.17%  '
$ws.Range("D16").Value = '58.002.47'
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").Value = '6.22'
$ws.Range("E17").Value = '  +10.23%  '
$ws.Range("D18").Value = '3.047.88'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").Value = '13.03'
$ws.Range("E19").Value = '  +5.91%  '
$ws.Range("E20").Value = '  +4.57%  '
$ws.Range("D21").Value = '338.05'
$ws.Range("E21").Value = '  +4.59%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '5.76'
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("E24").Value = '  +6.81%  '
$ws.Range("E25").Value = '  +5.18%  '
$ws.Range("E26").Value = '  +4.21%  '
$ws.Range("D27").Value = '0.0₃0955'
$ws.Range("E27").Value = '  +7.68%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = "'6.90"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.94%  '
$ws.Range("D30").Value = '7.57'
$ws.Range("E30").Value = '  +12.00%  '
$ws.Range("E31").Value = '  +5.39%  '
$ws.Range("E32").Value = '  +3.61%  '
$ws.Range("D33").Value = '21.03'
$ws.Range("E33").Value = '  +3.07%  '
$ws.Range("D34").Value = '4.77'
$ws.Range("E34").Value = '  +7.47%  '
$ws.Range("D35").Value = '155.77'
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("E36").Value = '  +7.14%  '
$ws.Range("E37").Value = '  +1.91%  '
$ws.Range("D38").Value = "'25.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +9.59%  '
$ws.Range("E39").Value = '  +3.00%  '
$ws.Range("D40").Value = '3.085.31'
$ws.Range("E40").Value = '  +2.54%  '
$ws.Range("D41").Value = '37.79'
$ws.Range("E41").Value = '  +4.43%  '
$ws.Range("E42").Value = '  +9.96%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  +3.85%  '
$ws.Range("D45").Value = '2.326.46'
$ws.Range("E45").Value = '  +3.62%  '
$ws.Range("E46").Value = '  +4.52%  '
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("E48").Value = '  +5.64%  '
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D50").Value = "'19.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.87%  '
$ws.Range("D51").Value = '1.87'
$ws.Range("E51").Value = '  -2.86%  '
